# "updating the files and cleanup"
#
# 1. Insert a new worksheet "Google_Search_Test" between the two existing
#    sheets (login_Test_Case_01, login_Test_Case_02) and populate it with
#    two rows of data ("TextToSearch" / "LinkedIn").
# 2. Move the active-sheet / selection state: login_Test_Case_01 is no
#    longer the selected tab (its selection becomes C23); the new
#    Google_Search_Test sheet becomes the active tab with an entire-row
#    selection on row 3 (A3:XFD3, active cell A3).

$wb = $excel.ActiveWorkbook

# --- Insert the new sheet right after "login_Test_Case_01" -----------------
$firstSheet = $wb.Worksheets.Item("login_Test_Case_01")
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $firstSheet)
$newSheet.Name = "Google_Search_Test"
$newSheet.Range("A1").Value = "TextToSearch"
$newSheet.Range("A2").Value = "LinkedIn"

# --- Update the (no-longer-active) first sheet's stored selection ----------
$ws1 = $wb.Worksheets.Item("login_Test_Case_01")
$ws1.Range("C23").Select() | Out-Null

# --- Make the new sheet active and select the full third row ---------------
$newSheet.Activate()
$newSheet.Range("A3").EntireRow.Select() | Out-Null
